$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "b.md" file. It is now ready for handoff. ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-14 02:51:47"

# --- zh-cn sheet: row 3 ("b.md") gets a fresh handoff + a version-mismatch error. ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "Content Duplicate" stores True/False as text, not a real boolean. Assigning the
# literal word "False" through .Value gets auto-coerced to a boolean by Excel, so
# copy the text value from a cell that already holds the text "False" instead.
$wsZhCn.Range("O2").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-14 02:51:38"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/5a547b6028d2232708f78c034eb18248f90f8ab2/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/86d812db6ec8c50e46a329558ae01f154114e3a7/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet: row 3 ("b.md") gets a fresh handoff + a version-mismatch error. ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("O2").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-14 02:51:47"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/5a547b6028d2232708f78c034eb18248f90f8ab2/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/86d812db6ec8c50e46a329558ae01f154114e3a7/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
